# Weekly refresh: insert a new price entry as the new row 91, pushing the
# existing rows 91-196 down to 92-197 (and growing the used range to R197).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 91 - shifts rows 91:196 down to 92:197.
$ws.Rows.Item(91).Insert()

# Populate the new row 91 with the latest weekly observation.
$ws.Range("A91").Value = 7
$ws.Range("B91").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C91").Value = "Ñuble"
$ws.Range("D91").Value = 44671
$ws.Range("E91").Value = 16
$ws.Range("F91").Value = 100112017
$ws.Range("G91").Value = "Apio"
$ws.Range("H91").Value = "Americana (o)"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 100
$ws.Range("K91").Value = 8000
$ws.Range("L91").Value = 8500
$ws.Range("M91").Value = 8250
$ws.Range("N91").Value = "$/docena de matas"
$ws.Range("O91").Value = "Provincia del Elquí"
$ws.Range("P91").Value = 1375
$ws.Range("Q91").Value = 6
$ws.Range("R91").Value = "Hortaliza"
